# adding file copy support
# Insert a new first worksheet "Company_cidr" ahead of the existing
# DC_networks / Sites_networks / Custom_identifier / NOT_TO_USE sheets,
# carrying a single header cell (styled like the other sheets' headers)
# plus an explanatory cell comment, and register the new shared string.

$wb = $excel.ActiveWorkbook

# --- create the new sheet and move it to the front of the tab strip ---
$ws = $wb.Worksheets.Add()
$ws.Name = "Company_cidr"
$ws.Move($wb.Worksheets.Item(1))

# --- header cell content ---
$ws.Range("A1").Value = "Company_network_range"

# --- match the pink/bold/bordered header style used on the other sheets ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.ThemeColor = 2
$ws.Range("A1").Interior.Color = 11353842
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").Borders.Weight = 2

# --- column sizing so column A comfortably fits the header text ---
$ws.Columns("A:A").ColumnWidth = 22.14

# --- explanatory note on the header cell ---
$comment = $ws.Range("A1").AddComment("Entire company network range`ne.g: 10.0.0.0/8`n")

# --- page orientation, matching the other worksheets ---
$ws.PageSetup.Orientation = 1

# --- leave the cursor where the author left it on the new sheet ---
$ws.Range("H10").Select()
